# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Updates the "base de datos" table (rows 16-19) on Hoja1:
#  - Row 16 now holds ANGELICA MARIA GULFO BASTIDAS / periodo 1708, with a
#    new "Valor Mora" of 200000 and Salario Basico reset to 0.
#  - Row 17 keeps ANGELICA MARIA GULFO BASTIDAS / periodo 1707, Salario
#    Basico reset to 0.
#  - Row 18 (ARLES MALDONADO WILCHES / periodo 1708) is unchanged.
#  - Row 19 now holds ARLES MALDONADO WILCHES / periodo 1707 with the
#    original 42300 / 1057500 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "52718112"
$ws.Range("D16").Value = "ANGELICA MARIA GULFO BASTIDAS"
$ws.Range("E16").Value = "1708"
$ws.Range("F16").Value = 200000
$ws.Range("G16").Value = 0

$ws.Range("G17").Value = 0

$ws.Range("C19").Value = "79908689"
$ws.Range("D19").Value = "ARLES MALDONADO WILCHES"
$ws.Range("E19").Value = "1707"
$ws.Range("F19").Value = 42300
$ws.Range("G19").Value = 1057500
